$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DO")

# Row 2: SV-N2-03 / dev6/port0/line27  ->  PV-FU-04 / dev6/port0/line24
$ws.Range("A2").Value = "PV-FU-04"
$ws.Range("B2").Value = "dev6/port0/line24"

# Row 3: HS_CAMERA / dev6/port0/line30  ->  HS_CAMERA / dev6/port0/line27
$ws.Range("B3").Value = "dev6/port0/line27"

# Row 4: SV-N2-02 / dev6/port0/line24  ->  SV-N2-02 / dev6/port0/line25
$ws.Range("B4").Value = "dev6/port0/line25"

# Row 5: empty  ->  SV-N2-03 / dev6/port0/line26 / NC / Closed
$ws.Range("A5").Value = "SV-N2-03"
$ws.Range("B5").Value = "dev6/port0/line26"
$ws.Range("C5").Value = "NC"
$ws.Range("D5").Value = "Closed"

# Update selection to E11 to match the saved view state
$ws.Range("E11").Select()
